# fix: update for christmas data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value  = "Santa Claus finger family for kids"
$ws.Range("A9").Value  = "Christmas Freeze Dance songs"
$ws.Range("A10").Value = "The Lights On the Christmas Tree"
$ws.Range("A11").Value = "Wish You a Merry Christmas | Carols"
$ws.Range("A22").Value = "Christmas Carols songs"
$ws.Range("A26").Value = "Hark The Herald Angel Sings"
$ws.Range("A37").Value = "Finger Family Santa Claus"
$ws.Range("A49").Value = "Joy To The World song"
$ws.Range("A51").Value = "Hello, Reindeer | Children's Christmas"
